$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.407949209213257
$ws.Range("B1").Value = 2.694223403930664
$ws.Range("C1").Value = 3.340522527694702
$ws.Range("D1").Value = 3.173795461654663
$ws.Range("E1").Value = 2.365808963775635
